$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (dates as Excel serial numbers, matching existing column A format)
$data = @(
    @(46070, 2110.53, 115398.8845837402, 115398.8845837402, 6843.22021484375, 1.184988498687744),
    @(46071, 2110.53, 116487.83,         116487.83,         6881.31005859375, 1.185185194015503),
    @(46072, 2110.53, 115796.91,         115796.91,         6861.89013671875, 1.178883910179138),
    @(46073, 2110.53, 116577.95,         116577.95,         6909.509765625,   1.176913619041443)
)

$startRow = 226
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

# Copy formatting (style/number format) from the last existing data row (225) down
# to the newly added rows, so column A keeps the date style (s="2").
$ws.Range("A225:F225").Copy()
$ws.Range("A226:F229").PasteSpecial(-4122)

[void]$ws.Range("A1").Select()
